{"js": "// Change: \"Pause: P (pressing the onscreen button or P to un-pause)\"\n//      -> \"Pause: P (pressing the onscreen button to un-pause)\"\n// i.e. remove the redundant \"or P \" so the only way to un-pause is pressing\n// the onscreen button (this is what the commit message \"made game easier to\n// beat\" is describing \u2014 fewer listed ways to resume the game).\n\nconst searchText = \"or P to un\";\nconst replaceText = \"to un\";\n\nconst results = context.document.body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(`Could not find text to edit: \"${searchText}\"`);\n}\n\n// Replace every match (there should be exactly one) in place, preserving\n// the run's existing character formatting.\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Change: \"Pause: P (pressing the onscreen button or P to un-pause)\"\n#      -> \"Pause: P (pressing the onscreen button to un-pause)\"\n# Removes the redundant \"or P \" so the bullet only calls out the onscreen\n# button as the way to un-pause (per commit message: \"made game easier to\n# beat\").\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"or P to un\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"to un\"\n\n# Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n#         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n# wdFindContinue = 1, wdReplaceAll = 2\n$found = $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"edit.ps1: could not find the text 'or P to un' to replace.\"\n}\n"}
